# R class gitkraken merge -> split "gitkraken" off with spell-check proofErr
# markers, add a blank line, and start a new "Change " paragraph that takes
# over the _GoBack bookmark (mirrors: place cursor at end of line 1, type a
# trailing space, press Enter twice, type "Change ").

$d = $word.ActiveDocument

function Escape-Xml {
    param([string]$s)
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

# --- Read the current state of paragraph 1 so the split point is computed,
# --- not hard-coded twice.
$firstPara = $d.Paragraphs(1).Range
$origText = $firstPara.Text
$needle = "gitkraken"
$idx = $origText.IndexOf($needle)
$beforeWord = $origText.Substring(0, $idx)
$afterWord = $origText.Substring($idx + $needle.Length)

$beforeXml = Escape-Xml $beforeWord
$afterXml = Escape-Xml $afterWord
$needleXml = Escape-Xml $needle

# --- Rebuild paragraph 1 as separate runs with spellcheck proofErr markers
# --- bracketing the non-dictionary word, plus the trailing space that was
# --- typed after "merge".
$para1Xml = @"
<w:p><w:r><w:t xml:space="preserve">$beforeXml</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>$needleXml</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">$afterXml</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
"@

# --- Two blank paragraphs (pressing Enter twice), then the new "Change "
# --- paragraph, which is where the cursor (and so _GoBack) ends up.
$para4Xml = '<w:p><w:r><w:t xml:space="preserve">Change </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$snippet = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>$para1Xml<w:p/><w:p/>$para4Xml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$d.Content.InsertXML($snippet)
